# Apply "Synced RM_Times" updates to the Event_Timing sheet.
# Updates Onset QRS1/Onset P/Onset QRS2/EMC1 (G,H,I,J) and the mirrored
# A/EMC2 columns (O,P) plus a couple of recomputed MD(ms) values (W)
# for several patient rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> list of (column letter, new value)
$updates = @{
    4  = @{ G = 111; H = 1146; I = 1013; J = 111; O = 1013; P = 1146 }
    5  = @{ H = 1066; I = 968;  O = 968;  P = 1066; W = 87.45 }
    7  = @{ H = 752;  I = 610;  O = 610;  P = 752 }
    11 = @{ G = 52;   H = 758;  J = 52;   P = 758 }
    13 = @{ H = 1060; I = 937;  O = 937;  P = 1060 }
    20 = @{ G = 116;  H = 1216; I = 1110; J = 116; O = 1110; P = 1216 }
    24 = @{ G = 98;   H = 1174; I = 1034; J = 98;  O = 1034; P = 1174 }
    25 = @{ G = 80;   H = 1281; I = 1174; J = 80;  O = 1174; P = 1281 }
    28 = @{ G = 58;   H = 974;  I = 890;  J = 58;  O = 890;  P = 974 }
    29 = @{ G = 10;   H = 959;  I = 879;  J = 10;  O = 879;  P = 959 }
    42 = @{ G = 46;   H = 635;  I = 546;  J = 46;  O = 546;  P = 635; W = 35.29 }
    53 = @{ H = 754;  I = 641;  O = 641;  P = 754; W = 73.64 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
